# Updates "cryptos" price/volume table in Sheet1 to reflect the latest
# GitHub Actions scrape (see commit message: "Updated cryptos list ...").
#
# Price values in column D are written with a leading single-quote so that
# Excel keeps them as plain text (matching the original inlineStr cells)
# instead of silently reinterpreting number-like strings (e.g. "0.9969",
# "334.22") as numeric values. Values that already contain two "." and are
# not parsed as numbers by Excel (e.g. "29.573.22") are written directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.573.22"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.920.17"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("D4").Value = "'0.9969"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "'334.22"
$ws.Range("E5").Value = "  -2.18%  "
$ws.Range("D6").Value = "'0.9971"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "'0.4648"
$ws.Range("E7").Value = "  -2.92%  "
$ws.Range("D8").Value = "'0.4156"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").Value = "'48.21"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "'0.08060"
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("D12").Value = "'22.40"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").Value = "1.915.10"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "'6.004"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").Value = "'7.176"
$ws.Range("E15").Value = "  -2.97%  "
$ws.Range("D16").Value = "'89.63"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "'0.9973"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "'0.00001036"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").Value = "'0.06592"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").Value = "'17.81"
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").Value = "'0.9992"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "29.530.89"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "'5.538"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "'11.47"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("D25").Value = "'2.199"
$ws.Range("E25").Value = "  -3.73%  "
$ws.Range("D26").Value = "2.134.40"
$ws.Range("E26").Value = "  -2.24%  "
$ws.Range("D27").Value = "'156.70"
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("D28").Value = "'19.93"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("D29").Value = "'2.166"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "'5.679"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "'117.56"
$ws.Range("E31").Value = "  -4.34%  "
$ws.Range("D32").Value = "'1.045"
$ws.Range("D33").Value = "'0.09463"
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").Value = "'1.442"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("D35").Value = "'5.451"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("D36").Value = "'3.536"
$ws.Range("E36").Value = "  -3.91%  "
$ws.Range("D37").Value = "'0.06138"
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("E38").Value = "  -2.16%  "
$ws.Range("D39").Value = "'8.479"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "'1.182"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").Value = "'0.5920"
$ws.Range("E41").Value = "  -2.79%  "
$ws.Range("D42").Value = "'0.9975"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").Value = "'10.28"
$ws.Range("E43").Value = "  -4.38%  "
$ws.Range("D44").Value = "'0.1840"
$ws.Range("E44").Value = "  -2.96%  "
$ws.Range("D45").Value = "'2.386"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("D46").Value = "'1.242"
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").Value = "'0.07544"
$ws.Range("E47").Value = "  +1.71%  "

# Rows 48 and 49: coin name/link swapped (Decentraland <-> EnergySwap)
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.5593"
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'12.21"
$ws.Range("E49").Value = "  -2.41%  "

$ws.Range("D50").Value = "'1.939"
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("D51").Value = "'112.82"
$ws.Range("E51").Value = "  -0.28%  "
